$d = $word.ActiveDocument

# 1. Delete the entire "Desfibrilador biventricular..." paragraph
#    (including its paragraph mark), merging the heading paragraph
#    directly with the following "Materiais:" paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Desfibrilador biventricular para terapia de ressincronização.*") {
        $p.Range.Delete()
        break
    }
}

# 2. Update the title/heading text.
$d.Content.Find.Execute(
    "Implante de CDI Biventricular (CRT-D Amplia" + [char]0x2122 + ")",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "CRT-D Amplia", 2)

# 3. Update the "Materiais" list items: add a leading bullet character
#    and simplify/normalize wording & codes.
$bullet = [char]0x2022

$d.Content.Find.Execute(
    "Gerador " + [char]0x2013 + " Amplia" + [char]0x2122 + " CRT-D",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "$bullet Gerador Amplia", 2)

$d.Content.Find.Execute(
    "Eletrodo Atrial " + [char]0x2013 + " 5076-52",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "$bullet Eletrodo 5076-52", 2)

$d.Content.Find.Execute(
    "Eletrodo Ventricular Esquerdo " + [char]0x2013 + " 4298/4299",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "$bullet Eletrodo VE 4298/4299", 2)

$d.Content.Find.Execute(
    "Eletrodo de Choque " + [char]0x2013 + " 6935M-62",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "$bullet Eletrodo 6935M62", 2)

$d.Content.Find.Execute(
    "Bainha " + [char]0x2013 + " 6250VIC",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "$bullet Bainha 6250VIC", 2)

$d.Content.Find.Execute(
    "Ferramenta de Corte",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "$bullet Ferramenta de corte", 2)

$d.Content.Find.Execute(
    "Guia 0.014",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "$bullet Guia 014", 2)

$d.Content.Find.Execute(
    "Subseletora",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "$bullet Subseletora", 2)

$d.Content.Find.Execute(
    "Introdutor " + [char]0x2013 + " 3",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "$bullet Introdutor " + [char]0x2013 + " 3", 2)
